$d = $word.ActiveDocument

# Locate "Supplementary Table 1." within the run that currently reads
# "Supplementary Table 1. Results from a pair of hemagglutination inhibition ("
# and make that portion of text bold, which causes Word to split the run
# into a bold run ("Supplementary Table 1.") and a non-bold run
# (" Results from a pair of hemagglutination inhibition (").
$rng = $d.Content
$found = $rng.Find.Execute("Supplementary Table 1.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Font.Bold = 1
}
